# Applies the 2023-05-14 cryptocurrency price-feed refresh to the sheet.
# Source data: cryptos.xlsx (row-by-row Price/Volume(1h) updates, plus the
# reordering of the Toncoin/Monero and MXToken/TheSandbox row pairs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry is the cell to update, its new plain-text value, and whether
# that text looks like a number (so Excel would otherwise auto-convert it,
# e.g. "1.033" -> 1.033). The original cells store these as literal text,
# so numeric-looking values are written with a temporary Text number
# format and then returned to the default "Normal" style.
$updates = @(
    @{Cell='D2'; Value='27.463.14'; Numeric=$false},
    @{Cell='E2'; Value='  +2.21%  '; Numeric=$false},
    @{Cell='D3'; Value='1.847.20'; Numeric=$false},
    @{Cell='E3'; Value='  +1.84%  '; Numeric=$false},
    @{Cell='D4'; Value='1.033'; Numeric=$true},
    @{Cell='E4'; Value='  +2.72%  '; Numeric=$false},
    @{Cell='D5'; Value='320.84'; Numeric=$true},
    @{Cell='E5'; Value='  +3.42%  '; Numeric=$false},
    @{Cell='D6'; Value='1.029'; Numeric=$true},
    @{Cell='E6'; Value='  +2.45%  '; Numeric=$false},
    @{Cell='D7'; Value='0.4362'; Numeric=$true},
    @{Cell='E7'; Value='  +1.57%  '; Numeric=$false},
    @{Cell='D8'; Value='0.3742'; Numeric=$true},
    @{Cell='E8'; Value='  +1.34%  '; Numeric=$false},
    @{Cell='D9'; Value='0.07373'; Numeric=$true},
    @{Cell='E9'; Value='  +1.87%  '; Numeric=$false},
    @{Cell='D10'; Value='0.8722'; Numeric=$true},
    @{Cell='E10'; Value='  +1.06%  '; Numeric=$false},
    @{Cell='D11'; Value='21.30'; Numeric=$true},
    @{Cell='E11'; Value='  +1.65%  '; Numeric=$false},
    @{Cell='D12'; Value='1.866.14'; Numeric=$false},
    @{Cell='E12'; Value='  -11.34%  '; Numeric=$false},
    @{Cell='D13'; Value='5.503'; Numeric=$true},
    @{Cell='E13'; Value='  +2.71%  '; Numeric=$false},
    @{Cell='D14'; Value='6.659'; Numeric=$true},
    @{Cell='E14'; Value='  +0.77%  '; Numeric=$false},
    @{Cell='D15'; Value='0.07211'; Numeric=$true},
    @{Cell='E15'; Value='  +4.19%  '; Numeric=$false},
    @{Cell='D16'; Value='82.54'; Numeric=$true},
    @{Cell='E16'; Value='  +2.14%  '; Numeric=$false},
    @{Cell='D17'; Value='1.035'; Numeric=$true},
    @{Cell='E17'; Value='  +2.51%  '; Numeric=$false},
    @{Cell='D18'; Value='0.000008993'; Numeric=$true},
    @{Cell='E18'; Value='  +1.04%  '; Numeric=$false},
    @{Cell='D19'; Value='1.029'; Numeric=$true},
    @{Cell='E19'; Value='  +2.40%  '; Numeric=$false},
    @{Cell='D20'; Value='15.38'; Numeric=$true},
    @{Cell='E20'; Value='  +1.17%  '; Numeric=$false},
    @{Cell='D21'; Value='27.490.97'; Numeric=$false},
    @{Cell='E21'; Value='  +2.09%  '; Numeric=$false},
    @{Cell='D22'; Value='5.240'; Numeric=$true},
    @{Cell='E22'; Value='  +1.11%  '; Numeric=$false},
    @{Cell='D23'; Value='11.18'; Numeric=$true},
    @{Cell='E23'; Value='  +0.77%  '; Numeric=$false},
    @{Cell='D24'; Value='2.083.46'; Numeric=$false},
    @{Cell='E24'; Value='  -10.94%  '; Numeric=$false},
    @{Cell='B25'; Value='Monero'; Numeric=$false},
    @{Cell='C25'; Value='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'; Numeric=$false},
    @{Cell='D25'; Value='157.45'; Numeric=$true},
    @{Cell='E25'; Value='  +2.32%  '; Numeric=$false},
    @{Cell='B26'; Value='Toncoin'; Numeric=$false},
    @{Cell='C26'; Value='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; Numeric=$false},
    @{Cell='D26'; Value='1.931'; Numeric=$true},
    @{Cell='E26'; Value='  +2.61%  '; Numeric=$false},
    @{Cell='D27'; Value='18.65'; Numeric=$true},
    @{Cell='E27'; Value='  +2.07%  '; Numeric=$false},
    @{Cell='D28'; Value='5.249'; Numeric=$true},
    @{Cell='E28'; Value='  +1.07%  '; Numeric=$false},
    @{Cell='D29'; Value='1.937'; Numeric=$true},
    @{Cell='E29'; Value='  +2.62%  '; Numeric=$false},
    @{Cell='D30'; Value='116.71'; Numeric=$true},
    @{Cell='E30'; Value='  +1.87%  '; Numeric=$false},
    @{Cell='D31'; Value='0.09029'; Numeric=$true},
    @{Cell='E31'; Value='  +1.05%  '; Numeric=$false},
    @{Cell='D32'; Value='1.195'; Numeric=$true},
    @{Cell='E32'; Value='  +2.41%  '; Numeric=$false},
    @{Cell='D33'; Value='0.7593'; Numeric=$true},
    @{Cell='E33'; Value='  +1.77%  '; Numeric=$false},
    @{Cell='D34'; Value='4.495'; Numeric=$true},
    @{Cell='E34'; Value='  +1.98%  '; Numeric=$false},
    @{Cell='D35'; Value='2.888'; Numeric=$true},
    @{Cell='E35'; Value='  +3.14%  '; Numeric=$false},
    @{Cell='D36'; Value='1.031'; Numeric=$true},
    @{Cell='E36'; Value='  +1.94%  '; Numeric=$false},
    @{Cell='D37'; Value='1.146'; Numeric=$true},
    @{Cell='E37'; Value='  +2.61%  '; Numeric=$false},
    @{Cell='D38'; Value='0.01967'; Numeric=$true},
    @{Cell='E38'; Value='  +2.44%  '; Numeric=$false},
    @{Cell='D39'; Value='0.05263'; Numeric=$true},
    @{Cell='E39'; Value='  +1.31%  '; Numeric=$false},
    @{Cell='B40'; Value='TheSandbox'; Numeric=$false},
    @{Cell='C40'; Value='https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'; Numeric=$false},
    @{Cell='D40'; Value='0.5144'; Numeric=$true},
    @{Cell='E40'; Value='  +1.55%  '; Numeric=$false},
    @{Cell='B41'; Value='MXToken'; Numeric=$false},
    @{Cell='C41'; Value='https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'; Numeric=$false},
    @{Cell='D41'; Value='2.803'; Numeric=$true},
    @{Cell='E41'; Value='  +3.43%  '; Numeric=$false},
    @{Cell='D42'; Value='0.1665'; Numeric=$true},
    @{Cell='E42'; Value='  +1.01%  '; Numeric=$false},
    @{Cell='D43'; Value='6.681'; Numeric=$true},
    @{Cell='E43'; Value='  +2.49%  '; Numeric=$false},
    @{Cell='D44'; Value='8.510'; Numeric=$true},
    @{Cell='E44'; Value='  +2.80%  '; Numeric=$false},
    @{Cell='D45'; Value='108.68'; Numeric=$true},
    @{Cell='E45'; Value='  +1.97%  '; Numeric=$false},
    @{Cell='D46'; Value='10.47'; Numeric=$true},
    @{Cell='E46'; Value='  +0.07%  '; Numeric=$false},
    @{Cell='D47'; Value='1.705'; Numeric=$true},
    @{Cell='E47'; Value='  +3.52%  '; Numeric=$false},
    @{Cell='D48'; Value='0.4630'; Numeric=$true},
    @{Cell='E48'; Value='  +1.94%  '; Numeric=$false},
    @{Cell='D49'; Value='0.06375'; Numeric=$true},
    @{Cell='E49'; Value='  +1.13%  '; Numeric=$false},
    @{Cell='D50'; Value='1.858'; Numeric=$true},
    @{Cell='E50'; Value='  +3.57%  '; Numeric=$false},
    @{Cell='D51'; Value='39.02'; Numeric=$true},
    @{Cell='E51'; Value='  +3.68%  '; Numeric=$false}
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)

    if ($u.Numeric) {
        $cell.NumberFormat = '@'
        $cell.Value = $u.Value
        $cell.Style = 'Normal'
    } else {
        $cell.Value = $u.Value
    }
}
